$wb = $excel.ActiveWorkbook

# Update the "Date" metadata value on the "Metadata" sheet
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2023-07-27T13:16:51+00:00"

# Rename the "prolongee-prorogee" code/display to "prorogee" / "Prorogée" on the "Concepts" sheet
$wsConcepts = $wb.Worksheets.Item("Concepts")
$wsConcepts.Range("B5").Value = "prorogee"
$wsConcepts.Range("C5").Value = "Prorogée"
